# Auto update ESO load data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rolling-window ESO load data for rows 2-7 (columns A..AA:
# Year, Month, Day, then 24 half-hourly/hourly settlement periods).
$rows = @{
    2 = @(2026, 2, 1, 5087, 4931, 4880, 4832, 4779, 4848, 4946, 5219, 5604, 6022, 6210, 6189, 6151, 6081, 6026, 6037, 6209, 6546, 6763, 6702, 6539, 6294, 6121, 5841)
    3 = @(2026, 2, 2, 5588, 5362, 5266, 5284, 5336, 5466, 5894, 6414, 6892, 7121, 7085, 6967, 6823, 6799, 6744, 6727, 6665, 6939, 7331, 7300, 7103, 6836, 6601, 6240)
    4 = @(2026, 2, 3, 5886, 5634, 5557, 5555, 5590, 5740, 6227, 6780, 7244, 7364, 7245, 7111, 6932, 6835, 6733, 6744, 6814, 7029, 7292, 7273, 7055, 6770, 6513, 6136)
    5 = @(2026, 2, 4, 5795, 5499, 5408, 5407, 5440, 5587, 6061, 6600, 7052, 7170, 7090, 7001, 6883, 6818, 6750, 6756, 6793, 6906, 7045, 7027, 6815, 6539, 6290, 5925)
    6 = @(2026, 2, 5, 5445, 5151, 5061, 5059, 5093, 5238, 5710, 6246, 6697, 6813, 6753, 6686, 6596, 6547, 6495, 6500, 6525, 6604, 6701, 6682, 6467, 6185, 5932, 5560)
    7 = @(2026, 2, 6, 5075, 4783, 4694, 4693, 4726, 4870, 5340, 5873, 6322, 6438, 6398, 6353, 6293, 6261, 6226, 6227, 6231, 6243, 6257, 6238, 6025, 5746, 5496, 5127)
}

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}
